$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column M header cell (row 2, thin/thick bottom border row) ---------
$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial(-4122)   # xlPasteFormats

# --- Column M year header (row 3) ---------------------------------------
$ws.Range("F3").Copy()
$ws.Range("M3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("M3").Value = 2021

# --- Data rows 4 & 5: drop the custom "0.0" number format back to the ---
# --- sheet's default (General) formatting while keeping font/align -----
$ws.Range("B4").Copy()
$ws.Range("D4:M5").PasteSpecial(-4122)

# --- Data row 6 (bottom-bordered) ---------------------------------------
$ws.Range("B6").Copy()
$ws.Range("D6:M6").PasteSpecial(-4122)

# --- Restore / set the numeric values (PasteSpecial formats only, so ----
# --- the pre-existing values in D:L survive untouched automatically; ----
# --- only the brand-new column M needs values) ---------------------------
$ws.Range("M4").Value = 7105
$ws.Range("M5").Value = 81079
$ws.Range("M6").Value = 214139

$excel.CutCopyMode = $false

Write-Output "done"
